$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting already used in column C (e.g. C11) down into
# C12:C13 so they share the same style record as the rest of the column.
$ws.Range("C11").Copy()
$ws.Range("C12:C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 12: McMaster purchase of screws/nuts
$ws.Range("A12").Value = "McMaster"
$ws.Range("B12").Value = 20.02
$ws.Range("C12").Value = 44785
$ws.Range("D12").Value = "#3-48 Screws and Nuts"

# Row 13: PJRC purchase of Teensy 4.1s
$ws.Range("A13").Value = "PJRC"
$ws.Range("B13").Value = 97.18
$ws.Range("C13").Value = 44785
$ws.Range("D13").Value = "Teensy 4.1s"

# Update selection to A14
$ws.Range("A14").Select()

# Column C now needs an explicit best-fit width, since the new dates no
# longer fit the default column width.
$ws.Columns.Item(3).EntireColumn.AutoFit()
